$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Routes")

# Update row 4 to reflect the new route: Mexico City (MMMX) -> Seattle (KSEA)
$ws.Range("C4").Value = "MMMX"
$ws.Range("B4").Value = "Aeropuerto México Ciudad Intl"
$ws.Range("D4").Value = "Seattle Tacoma Intl"
$ws.Range("E4").Value = "KSEA"

$ws.Range("B5").Select()
